# This script re-applies the latest (per-hour) snapshot of the cryptos
# table: columns B (Coin), C (Link), D (Price) and E (Volume(1h)) for
# rows 2-51 of the active sheet are refreshed to the new scraped values.
# A couple of rows also swap rank position with their neighbour (e.g.
# Litecoin/Dai, USDe/Fetch.AI, Mantle/EnergySwap, Maker/Hedera/OKB) so
# both the B/C (name/link) and D/E (price/volume) cells are rewritten
# for those rows.
#
# Note: some Price values are plain decimal-looking text (e.g.
# "582.59"). Assigning such text to Range.Value makes Excel treat it
# like manual keyboard entry and auto-convert it to a real number
# (losing the original text formatting used throughout this sheet).
# To keep parity with the source data -- which stores every Price/
# Volume cell as text -- those assignments are prefixed with a leading
# apostrophe, Excel's standard "treat as text" quote-prefix, exactly
# as if a user had typed `'582.59` into the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.003.38'
$ws.Range('E2').Value = '  +0.30%  '

$ws.Range('D3').Value = '3.245.20'
$ws.Range('E3').Value = '  -0.31%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').Value = '''582.59'
$ws.Range('E5').Value = '  +0.65%  '

$ws.Range('D6').Value = '''185.66'
$ws.Range('E6').Value = '  +2.11%  '

$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('D8').Value = '''0.595'
$ws.Range('E8').Value = '  -0.42%  '

$ws.Range('D9').Value = '''0.131'
$ws.Range('E9').Value = '  -1.88%  '

$ws.Range('D10').Value = '''6.66'
$ws.Range('E10').Value = '  -0.12%  '

$ws.Range('D11').Value = '''0.419'
$ws.Range('E11').Value = '  +0.82%  '

$ws.Range('D12').Value = '3.808.93'
$ws.Range('E12').Value = '  -0.25%  '

$ws.Range('D13').Value = '''0.138'
$ws.Range('E13').Value = '  +0.12%  '

$ws.Range('D14').Value = '''28.07'
$ws.Range('E14').Value = '  -1.86%  '

$ws.Range('D15').Value = '68.057.25'
$ws.Range('E15').Value = '  +0.42%  '

$ws.Range('D16').Value = '''0.0000171'
$ws.Range('E16').Value = '  -0.70%  '

$ws.Range('D17').Value = '3.215.52'
$ws.Range('E17').Value = '  -1.19%  '

$ws.Range('D18').Value = '''5.82'
$ws.Range('E18').Value = '  -0.31%  '

$ws.Range('D19').Value = '''13.55'
$ws.Range('E19').Value = '  +0.17%  '

$ws.Range('D20').Value = '''395.13'
$ws.Range('E20').Value = '  +4.19%  '

$ws.Range('D21').Value = '''7.62'
$ws.Range('E21').Value = '  -0.40%  '

$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '''1.00'
$ws.Range('E22').Value = '  +0.05%  '

$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').Value = '''71.48'
$ws.Range('E23').Value = '  +0.11%  '

$ws.Range('E24').Value = '  +1.06%  '

$ws.Range('E25').Value = '  -0.40%  '

$ws.Range('E26').Value = '  +3.60%  '

$ws.Range('D27').Value = '''9.77'
$ws.Range('E27').Value = '  -1.73%  '

$ws.Range('E28').Value = '  -0.08%  '

$ws.Range('E29').Value = '  -0.50%  '

$ws.Range('D30').Value = '''5.65'
$ws.Range('E30').Value = '  -0.18%  '

$ws.Range('D31').Value = '''22.83'
$ws.Range('E31').Value = '  -0.46%  '

$ws.Range('D32').Value = '''7.12'
$ws.Range('E32').Value = '  +1.74%  '

$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Value = '''1.27'
$ws.Range('E33').Value = '  +0.20%  '

$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').Value = '''1.00'
$ws.Range('E34').Value = '  +0.16%  '

$ws.Range('D35').Value = '''162.27'
$ws.Range('E35').Value = '  -0.83%  '

$ws.Range('D36').Value = '''1.50'
$ws.Range('E36').Value = '  -4.07%  '

$ws.Range('D37').Value = '''1.92'
$ws.Range('E37').Value = '  +1.59%  '

$ws.Range('B38').Value = 'Mantle'
$ws.Range('C38').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D38').Value = '''0.816'
$ws.Range('E38').Value = '  -3.78%  '

$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').Value = '''26.48'
$ws.Range('E39').Value = '  -0.22%  '

$ws.Range('E40').Value = '  -1.11%  '

$ws.Range('D41').Value = '''6.53'
$ws.Range('E41').Value = '  -1.80%  '

$ws.Range('E42').Value = '  -4.33%  '

$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').Value = '''0.0690'
$ws.Range('E43').Value = '  +1.02%  '

$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '''41.08'
$ws.Range('E44').Value = '  +0.07%  '

$ws.Range('D45').Value = '''25.18'
$ws.Range('E45').Value = '  -1.59%  '

$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '2.616.96'
$ws.Range('E46').Value = '  +0.17%  '

$ws.Range('D47').Value = '''338.59'
$ws.Range('E47').Value = '  -2.85%  '

$ws.Range('D48').Value = '''0.0280'
$ws.Range('E48').Value = '  -1.77%  '

$ws.Range('D49').Value = '''6.36'
$ws.Range('E49').Value = '  +3.11%  '

$ws.Range('E50').Value = '  -1.01%  '

$ws.Range('D51').Value = '''31.29'
$ws.Range('E51').Value = '  +2.43%  '
